function Set-TextValue($Worksheet, $CellRef, $Text) {
    $rng = $Worksheet.Range($CellRef)
    # Force text number-format first so numeric-looking strings (e.g. "1.0000",
    # "15.13") are stored as literal text rather than being parsed into numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    # Reset the style back to the workbook default so no stray cell-style diff is
    # left behind (matches the original file, where these cells carry no explicit s=).
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '25.895.61'
Set-TextValue $ws 'E2' '  +0.00%  '
Set-TextValue $ws 'D3' '1.730.85'
Set-TextValue $ws 'E3' '  -0.59%  '
Set-TextValue $ws 'E4' '  +0.03%  '
Set-TextValue $ws 'D5' '244.98'
Set-TextValue $ws 'E5' '  +2.77%  '
Set-TextValue $ws 'E6' '  +0.01%  '
Set-TextValue $ws 'D7' '0.5030'
Set-TextValue $ws 'E7' '  -2.81%  '
Set-TextValue $ws 'D8' '0.2708'
Set-TextValue $ws 'E8' '  -1.57%  '
Set-TextValue $ws 'D9' '0.06164'
Set-TextValue $ws 'E9' '  +0.09%  '
Set-TextValue $ws 'D10' '1.735.68'
Set-TextValue $ws 'E10' '  -0.31%  '
Set-TextValue $ws 'E11' '  +0.91%  '
Set-TextValue $ws 'B12' 'Solana'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws 'D12' '15.13'
Set-TextValue $ws 'E12' '  +0.83%  '
Set-TextValue $ws 'B13' 'Polygon'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D13' '0.6490'
Set-TextValue $ws 'E13' '  +0.61%  '
Set-TextValue $ws 'D14' '4.749'
Set-TextValue $ws 'E14' '  +3.15%  '
Set-TextValue $ws 'D15' '76.89'
Set-TextValue $ws 'E15' '  -0.75%  '
Set-TextValue $ws 'D16' '1.0000'
Set-TextValue $ws 'E16' '  -0.06%  '
Set-TextValue $ws 'D17' '1.001'
Set-TextValue $ws 'E17' '  +0.13%  '
Set-TextValue $ws 'D18' '25.883.07'
Set-TextValue $ws 'E18' '  -0.12%  '
Set-TextValue $ws 'E19' '  +1.40%  '
Set-TextValue $ws 'D20' '0.000006819'
Set-TextValue $ws 'E20' '  +0.55%  '
Set-TextValue $ws 'D21' '4.588'
Set-TextValue $ws 'E21' '  +7.11%  '
Set-TextValue $ws 'D22' '1.955.20'
Set-TextValue $ws 'E22' '  -0.53%  '
Set-TextValue $ws 'D23' '8.781'
Set-TextValue $ws 'E23' '  +1.44%  '
Set-TextValue $ws 'D24' '5.477'
Set-TextValue $ws 'E24' '  +4.11%  '
Set-TextValue $ws 'D25' '134.50'
Set-TextValue $ws 'E25' '  -3.41%  '
Set-TextValue $ws 'D26' '15.24'
Set-TextValue $ws 'E26' '  +0.74%  '
Set-TextValue $ws 'D27' '1.419'
Set-TextValue $ws 'E27' '  -6.12%  '
Set-TextValue $ws 'D28' '1.782'
Set-TextValue $ws 'E28' '  +1.12%  '
Set-TextValue $ws 'D29' '105.12'
Set-TextValue $ws 'E29' '  -0.72%  '
Set-TextValue $ws 'D30' '3.957'
Set-TextValue $ws 'E30' '  +0.61%  '
Set-TextValue $ws 'E31' '  -1.98%  '
Set-TextValue $ws 'D32' '3.695'
Set-TextValue $ws 'E32' '  -0.03%  '
Set-TextValue $ws 'D33' '0.04719'
Set-TextValue $ws 'E33' '  +2.70%  '
Set-TextValue $ws 'D34' '2.656'
Set-TextValue $ws 'E34' '  +0.56%  '
Set-TextValue $ws 'D35' '0.9950'
Set-TextValue $ws 'E35' '  +0.76%  '
Set-TextValue $ws 'D36' '0.6107'
Set-TextValue $ws 'E36' '  -1.28%  '
Set-TextValue $ws 'D37' '2.743'
Set-TextValue $ws 'E37' '  +2.23%  '
Set-TextValue $ws 'D38' '0.01606'
Set-TextValue $ws 'E38' '  -0.17%  '
Set-TextValue $ws 'D39' '0.8708'
Set-TextValue $ws 'E39' '  +17.47%  '
Set-TextValue $ws 'D40' '1.950'
Set-TextValue $ws 'E40' '  +1.13%  '
Set-TextValue $ws 'D41' '1.001'
Set-TextValue $ws 'E41' '  +0.05%  '
Set-TextValue $ws 'D42' '101.88'
Set-TextValue $ws 'E42' '  +4.07%  '
Set-TextValue $ws 'D43' '0.3882'
Set-TextValue $ws 'E43' '  +1.15%  '
Set-TextValue $ws 'D44' '5.004'
Set-TextValue $ws 'E44' '  +0.37%  '
Set-TextValue $ws 'D45' '0.1182'
Set-TextValue $ws 'E45' '  +4.64%  '
Set-TextValue $ws 'D46' '6.341'
Set-TextValue $ws 'E46' '  +2.06%  '
Set-TextValue $ws 'D47' '55.58'
Set-TextValue $ws 'E47' '  +1.12%  '
Set-TextValue $ws 'D48' '0.05279'
Set-TextValue $ws 'E48' '  +0.40%  '
Set-TextValue $ws 'D49' '30.74'
Set-TextValue $ws 'E49' '  +0.91%  '
Set-TextValue $ws 'B50' 'EnergySwap'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D50' '7.678'
Set-TextValue $ws 'E50' '  +0.95%  '
Set-TextValue $ws 'B51' 'Decentraland'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws 'D51' '0.3477'
Set-TextValue $ws 'E51' '  +2.08%  '
